$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.144.10"
$ws.Range("E2").Value = "  -1.84%  "

$ws.Range("D3").Value = "2.289.91"
$ws.Range("E3").Value = "  -3.16%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.57"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.03"
$ws.Range("E6").Value = "  -4.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  -1.31%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.603"
$ws.Range("E9").Value = "  -3.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.14"
$ws.Range("E10").Value = "  -5.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0906"
$ws.Range("E11").Value = "  -2.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.23"
$ws.Range("E12").Value = "  -4.93%  "

$ws.Range("E13").Value = "  -0.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.956"
$ws.Range("E14").Value = "  -5.73%  "

$ws.Range("E15").Value = "  -5.10%  "

$ws.Range("D16").Value = "2.636.87"
$ws.Range("E16").Value = "  -3.17%  "

$ws.Range("D17").Value = "2.296.74"
$ws.Range("E17").Value = "  -3.47%  "

$ws.Range("D18").Value = "41.912.29"
$ws.Range("E18").Value = "  -2.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.33"
$ws.Range("E19").Value = "  -3.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000105"
$ws.Range("E20").Value = "  -1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.20"
$ws.Range("E21").Value = "  -4.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.60"
$ws.Range("E22").Value = "  -1.44%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.25"
$ws.Range("E23").Value = "  +3.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.20"
$ws.Range("E24").Value = "  +7.72%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("E25").Value = "  -3.04%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.43"
$ws.Range("E27").Value = "  +7.43%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.76"
$ws.Range("E28").Value = "  -6.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.71"
$ws.Range("E29").Value = "  -3.26%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.95"
$ws.Range("E30").Value = "  -2.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "163.48"
$ws.Range("E31").Value = "  -2.83%  "

$ws.Range("E32").Value = "  -4.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.82"
$ws.Range("E33").Value = "  -3.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.76"
$ws.Range("E34").Value = "  -4.32%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.135"
$ws.Range("E35").Value = "  +2.61%  "

$ws.Range("E36").Value = "  -5.61%  "

$ws.Range("E37").Value = "  -5.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0345"
$ws.Range("E38").Value = "  -5.21%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.81"
$ws.Range("E39").Value = "  +2.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.73"
$ws.Range("E40").Value = "  -3.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.27"
$ws.Range("E41").Value = "  -6.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.06"
$ws.Range("E43").Value = "  -3.32%  "

$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("E45").Value = "  -6.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.88"
$ws.Range("E46").Value = "  -4.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "112.23"
$ws.Range("E47").Value = "  -1.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.54"
$ws.Range("E48").Value = "  -1.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.88"
$ws.Range("E49").Value = "  -3.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.26"
$ws.Range("E50").Value = "  -5.38%  "

$ws.Range("D51").Value = "1.587.84"
$ws.Range("E51").Value = "  +0.56%  "

